$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:G, rows 2-11 (Q0..Q9)
$data = @(
    @(0.248767871177931, 0.8568015655085476, 2.374591342905846, 1.540970909169231, 1.536197893679839, 50),
    @(0.04882215121843843, 0.7469394130861163, 1.457578478961404, 1.207302148992291, 1.219356202157649, 47),
    @(0.05541629173989827, 0.751460361221017, 1.308898174481467, 1.144070878259501, 1.155355179810744, 46),
    @(0.2921665240571426, 0.8607796199794472, 2.535407960140988, 1.592296442293642, 1.582558739570702, 46),
    @(0.07348626165398647, 0.8093869007750452, 1.94618247269454, 1.395056440684225, 1.409225572003233, 44),
    @(-0.009028182342375235, 0.9204870629536859, 2.016096326641489, 1.419893068734927, 1.442583638701003, 32),
    @(0.07140803266404031, 0.918135862785761, 1.745723223468771, 1.321258197124533, 1.341135683000723, 31),
    @(0.135685191380836, 0.7164848718710426, 0.9226493013901432, 0.9605463556695967, 0.9868108548879626, 14),
    @(0.3151042569182115, 0.8922117577591338, 1.375200986598577, 1.172689637797903, 1.207552563243096, 8),
    @(0.6609969171670846, 0.6609969171670846, 0.9723557505147084, 0.9860810060612203, 0.8449369412449811, 4)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j  # Column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
